$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 8, pushing existing rows 8-20 down to 9-21
$ws.Rows.Item(8).Insert()

# Populate the new row 8 with the new Chirimoya price record
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C8").Value = "Ñuble"
$ws.Range("D8").Value = 45189
$ws.Range("D8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E8").Value = 16
$ws.Range("F8").Value = "Fruta"
$ws.Range("G8").Value = 100107
$ws.Range("H8").Value = "Otros"
$ws.Range("I8").Value = 100107002
$ws.Range("J8").Value = "Chirimoya"
$ws.Range("K8").Value = "Cultivar IV Región"
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 50
$ws.Range("N8").Value = 22000
$ws.Range("O8").Value = 22000
$ws.Range("P8").Value = 22000
$ws.Range("Q8").Value = "$/bandeja 10 kilos"
$ws.Range("R8").Value = "Provincia de Limarí"
$ws.Range("S8").Value = 2200
$ws.Range("T8").Value = 10
